# Auto update stock data
# Roll the report date forward from 2026/01/10 to 2026/01/11 for every
# ticker's latest-data row, and refresh Jabil's EBITDA figure for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "2026/01/11"
}

# Jabil (row 62) EBITDA value changed from 11.30 to 11.25
$ebitda = $ws.Cells.Item(62, 2)
$ebitda.NumberFormat = "@"
$ebitda.Value = "11.25"
